$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 10:05"

# --- Row 6 (Rusia) - updated case counts ---
$ws.Cells.Item(6, 2).Value = 370680
$ws.Cells.Item(6, 3).Value = 8338
$ws.Cells.Item(6, 4).Value = 142208
$ws.Cells.Item(6, 5).Value = 224504
$ws.Cells.Item(6, 7).Value = 161
$ws.Cells.Item(6, 8).Value = 3968

# --- Rows 72-73: Sudan now sorts ahead of Luxemburgo ---
$ws.Range("A72").Value = "Sudan"
$ws.Cells.Item(72, 2).Value = 4146
$ws.Cells.Item(72, 3).Value = 170
$ws.Cells.Item(72, 4).Value = 588
$ws.Cells.Item(72, 5).Value = 3374
$ws.Cells.Item(72, 7).Value = 14
$ws.Cells.Item(72, 8).Value = 184

$ws.Range("A73").Value = "Luxemburgo"
$ws.Cells.Item(73, 2).Value = 3995
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 3783
$ws.Cells.Item(73, 5).Value = 102
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 110

# --- Rows 82-86: "Consejo Danes para los Refugiados" now sorts ahead of
#     Costa de Marfil / Republica de Yibuti / Bulgaria / Bosnia y Herzegovina
#     (each row's country shifts down one slot) ---
$ws.Range("A82").Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(82, 2).Value = 2546
$ws.Cells.Item(82, 3).Value = 143
$ws.Cells.Item(82, 4).Value = 365
$ws.Cells.Item(82, 5).Value = 2113
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 68

$ws.Range("A83").Value = "Costa de Marfil"
$ws.Cells.Item(83, 2).Value = 2477
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 1286
$ws.Cells.Item(83, 5).Value = 1161
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 30

$ws.Range("A84").Value = "Republica de Yibuti"
$ws.Cells.Item(84, 2).Value = 2468
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 1079
$ws.Cells.Item(84, 5).Value = 1375
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 14

$ws.Range("A85").Value = "Bulgaria"
$ws.Cells.Item(85, 2).Value = 2460
$ws.Cells.Item(85, 3).Value = 17
$ws.Cells.Item(85, 4).Value = 912
$ws.Cells.Item(85, 5).Value = 1415
$ws.Cells.Item(85, 7).Value = 3
$ws.Cells.Item(85, 8).Value = 133

$ws.Range("A86").Value = "Bosnia y Herzegovina"
$ws.Cells.Item(86, 2).Value = 2416
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 1721
$ws.Cells.Item(86, 5).Value = 546
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 149

# --- Row 92 (Estonia) - updated case counts ---
$ws.Cells.Item(92, 2).Value = 1840
$ws.Cells.Item(92, 3).Value = 6
$ws.Cells.Item(92, 4).Value = 1561
$ws.Cells.Item(92, 5).Value = 213
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 66

# --- Row 98 (Eslovaquia) - updated case counts ---
$ws.Cells.Item(98, 2).Value = 1515
$ws.Cells.Item(98, 3).Value = 2
$ws.Cells.Item(98, 4).Value = 1327
$ws.Cells.Item(98, 5).Value = 160

# --- Row 110 (Letonia) - updated case counts ---
$ws.Cells.Item(110, 2).Value = 1057
$ws.Cells.Item(110, 3).Value = 4
$ws.Cells.Item(110, 5).Value = 293
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 23

# --- Row 141 (Estado de Palestina) - updated case counts ---
$ws.Cells.Item(141, 2).Value = 434
$ws.Cells.Item(141, 3).Value = 5
$ws.Cells.Item(141, 5).Value = 66
